# Updated Tracking Sheet with Vidushi's Task Details
$wb = $excel.ActiveWorkbook

$defines = $wb.Worksheets.Item("Defines")
$tracking = $wb.Worksheets.Item("Tracking Sheet")

# ---------------------------------------------------------------------------
# 1. Correct team-member name on the Defines sheet (Vidhushi -> Vidushi) and
#    broaden the module name (GPS -> GPS/Compass/Bridge)
# ---------------------------------------------------------------------------
$defines.Range("B4").Value = "Vidushi"
$defines.Range("A4").Value = "GPS/Compass/Bridge"

# ---------------------------------------------------------------------------
# 2. Row height for the existing wrapped-description row (row 2 keeps the
#    sheet's default row height)
# ---------------------------------------------------------------------------
$tracking.Rows.Item(3).RowHeight = 27.6

# ---------------------------------------------------------------------------
# 3. Append Vidushi's three new timesheet entries
# ---------------------------------------------------------------------------

# --- Row 4: WIKI Report / Documentation -----------------------------------
$tracking.Range("A4").Value = 43553
$tracking.Range("C4").Value = "Documentation"
$tracking.Range("D4").Value = "Vidushi "
$tracking.Range("E4").Value = "1. Updated wiki schedule and BOM`n2. Updated High Level Software Architecture and Team Member details on Wiki. "
$tracking.Range("F4").Value = 6
$tracking.Range("G4").Value = "WIKI Report"
$tracking.Rows.Item(4).RowHeight = 41.4

# --- Row 5: BRIDGE / Interfacing, Coding & Testing -------------------------
$tracking.Range("A5").Value = 43554
$tracking.Range("C5").Value = "Interfacing, Coding & Testing"
$tracking.Range("D5").Value = "Vidushi "
$tracking.Range("E5").Value = "1. Interfacing of Bluetooth HC-05 Module with Sjone board and 5V power supply.`n2. Worked on Data Mode(Rx/Tx) of HC-05.`n3. Tested interfacing and working(Rx/TX) with Serial Monitor.`n"
$tracking.Range("F5").Value = 6
$tracking.Range("G5").Value = "BRIDGE"
$tracking.Rows.Item(5).RowHeight = 69

# --- Row 6: BRIDGE / Interfacing, Coding & Testing -------------------------
$tracking.Range("A6").Value = 43555
$tracking.Range("C6").Value = "Interfacing, Coding & Testing"
$tracking.Range("D6").Value = "Vidushi "
$tracking.Range("E6").Value = "1. Interfacing of Bluetooth HC-05 Module withUBS TTL Convertor to configure Module.`n2. Worked on Command Mode of HC-05 and configured with name Tech_Savy and changed config setting to 38400.`n3. Refactoring of code for Bridge communication`n4. Implemented C wrapper for UART2.cpp, uart_dev.cpp and for switch and LED files. `n5. Testing done with Sample HC-05 BLE Android Application."
$tracking.Range("F6").Value = "Almost Full Day"
$tracking.Range("G6").Value = "BRIDGE"
$tracking.Rows.Item(6).RowHeight = 110.4

# ---------------------------------------------------------------------------
# 4. Re-style every data row (2-6) in one pass: smaller font, top alignment,
#    wrap text on the description column, bold date column with the new
#    date format. Doing this as a single pass over the whole block (instead
#    of per-row) keeps the generated style table compact.
# ---------------------------------------------------------------------------
$allData = $tracking.Range("A2:G6")
$allData.Font.Size = 10
$allData.VerticalAlignment = -4160   # xlTop
$allData.HorizontalAlignment = -4131 # xlLeft

$allDateCells = $tracking.Range("A2:A6")
$allDateCells.Font.Bold = $true
$allDateCells.NumberFormat = "[$-409]d\-mmm\-yy;@"

$allDescCells = $tracking.Range("E3:E6")
$allDescCells.WrapText = $true

# ---------------------------------------------------------------------------
# 5. Restore selections to match the refreshed workbook view
# ---------------------------------------------------------------------------
$defines.Range("B14").Select()
$tracking.Activate()
$tracking.Range("D6").Select()
